# Apply the "accessibility panel" localisation update described by the commit:
#   - Add 12 new rows (keys + FR/EN labels) to the Labels sheet for the new
#     "accessibilityPanel.*" strings.
#   - Simplify the rich-text "CO2" (CO + subscript 2) label to plain "CO2".
#   - Move the active/selected tab from "Choices" to "Labels", with the
#     Labels selection moving from C20 to D20.

$wb = $excel.ActiveWorkbook

$labels = $wb.Worksheets.Item("Labels")
$choices = $wb.Worksheets.Item("Choices")

# ---------------------------------------------------------------------------
# 1. New localisation rows 34-45 on the "Labels" sheet.
#    Columns: A=section, B=path, C=label::fr, D=label::en
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row=34; A="results"; B="accessibilityPanel.title";             C="Accessibilité";      D="Accessibility" },
    @{ Row=35; A="results"; B="accessibilityPanel.locationsTitle";    C="Logement";           D="Locations" },
    @{ Row=36; A="results"; B="accessibilityPanel.bothAddresses";     C="Les deux";           D="Both" },
    @{ Row=37; A="results"; B="accessibilityPanel.firstAddressOnly";  C="Logement #1";        D="House #1" },
    @{ Row=38; A="results"; B="accessibilityPanel.secondAddressOnly"; C="Logement #2";        D="House #2" },
    @{ Row=39; A="results"; B="accessibilityPanel.travelTimeTitle";   C="Temps de trajet";    D="Travel time" },
    @{ Row=40; A="results"; B="accessibilityPanel.15min";             C="15 min.";            D="15 min." },
    @{ Row=41; A="results"; B="accessibilityPanel.30min";             C="30 min.";            D="30 min." },
    @{ Row=42; A="results"; B="accessibilityPanel.45min";             C="45 min.";            D="45 min." },
    @{ Row=43; A="results"; B="accessibilityPanel.modeOfTransportTitle"; C="Mode de transport"; D="Mode of transport" },
    @{ Row=44; A="results"; B="accessibilityPanel.minimize";          C="Réduire la fenêtre"; D="Minimize the panel" },
    @{ Row=45; A="results"; B="accessibilityPanel.expand";            C="Agrandir la fenêtre"; D="Maximize the panel" }
)

foreach ($r in $newRows) {
    $labels.Range("A" + $r.Row).Value = $r.A
    $labels.Range("B" + $r.Row).Value = $r.B
    $labels.Range("C" + $r.Row).Value = $r.C
    $labels.Range("D" + $r.Row).Value = $r.D
}

# ---------------------------------------------------------------------------
# 2. Simplify the "CO2" rich-text label (CO + subscript 2) down to plain text.
#    (locationComparison.environmentCo2 -> "CO2", row 25 on Labels)
#    C25/D25 both share the same rich-text string; temporarily move D25 off
#    of it so C25 becomes the sole owner, normalise C25's formatting in
#    place (drops the subscript run), then point D25 back at the now-plain
#    "CO2" text so both cells end up sharing the single simplified string.
# ---------------------------------------------------------------------------
$labels.Range("D25").Value = "__tmp_co2__"
$labels.Range("C25").Characters(1, 3).Font.Subscript = $false
$labels.Range("D25").Value = "CO2"

# ---------------------------------------------------------------------------
# 3. Switch the active/selected sheet from "Choices" to "Labels" and move the
#    Labels selection from C20 to D20.
# ---------------------------------------------------------------------------
$choices.Select()
$labels.Select()
$labels.Activate()
$labels.Range("D20").Select()
